$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the new student, shifting the Total row down
$ws.Rows.Item(5).Insert()

# New student row (row 5)
$ws.Cells.Item(5, 1).Value = "40724261 (Gastao Bettencourt)"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 10

# Code Review 3 column (D) values for existing students
$ws.Cells.Item(2, 4).Value = 30
$ws.Cells.Item(3, 4).Value = 30
$ws.Cells.Item(4, 4).Value = 30

# Total row (now row 6)
$ws.Cells.Item(6, 3).Value = 100
$ws.Range("D6").Formula = "=SUM(D2:D5)"

# Column A width (stored OOXML width ~27.54 -> nearest reachable pixel-grid value)
$ws.Columns.Item(1).ColumnWidth = 26.666666666666668

# Selection
$ws.Range("D10").Select()
